$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.577380418777466
$ws.Range("B1").Value = 1.82272744178772
$ws.Range("C1").Value = 2.307337045669556
$ws.Range("D1").Value = 3.748313188552856
$ws.Range("E1").Value = 3.067660808563232
